$wb = $excel.ActiveWorkbook

# --- "report" sheet: insert a new column before column C for the new MTS_msk
# (mts_msc) entry, copying date-format from the neighbouring column and
# filling in the five populated rows (customer_name, date, project_title,
# project_folder, supportsave_folder). -----------------------------------
$ws1 = $wb.Worksheets.Item("report")

$ws1.Columns("C:C").Insert()

# Match column B's width on the freshly inserted column C.
$ws1.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# Row 2 - customer_name: reuse the existing "MTS_msk" label (same value as
# the later "MTS_Techblock" column's neighbour before it).
$ws1.Range("C2").Value2 = "MTS_msk"

# Row 3 - date: copy number formatting from the adjoining date cell, then
# set the new date value (2022-01-14).
$ws1.Range("D3").Copy()
$ws1.Range("C3").PasteSpecial(-4122)
$ws1.Range("C3").Value2 = 44575

# Row 4 - project_title
$ws1.Range("C4").Value2 = "SAN Assessment"

# Row 5 - project_folder
$ws1.Range("C5").Value2 = "C:\Users\vlasenko\OneDrive - Hewlett Packard Enterprise\Documents\01.CUSTOMERS\MTS\SAN Assessment\JAN2022\mts_msc"

# Row 6 - supportsave_folder
$ws1.Range("C6").Value2 = "C:\Users\vlasenko\OneDrive - Hewlett Packard Enterprise\Documents\06.CONFIGS\MTS\JAN22\mts_msc\ssave"

$excel.CutCopyMode = $false

# --- view/selection bookkeeping (best effort) -----------------------------
$ws3 = $wb.Worksheets.Item("data_dependency")
$ws2 = $wb.Worksheets.Item("service_tables")

$ws1.Activate()
$ws1.Range("E22").Select()

$ws3.Activate()
$ws3.Range("D28").Select()

$ws2.Activate()
$ws2.Range("H20").Select()
